$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.111.96"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "1.934.07"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'242.07"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4794"
$ws.Range("D8").Value = "'0.2911"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "'0.06775"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'20.09"
$ws.Range("E10").Value = "  +6.94%  "
$ws.Range("D11").Value = "'104.37"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'0.07843"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "1.938.38"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'5.296"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "'0.6962"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "'295.89"
$ws.Range("E16").Value = "  +11.32%  "
$ws.Range("D17").Value = "31.110.79"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").Value = "'12.98"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.195.21"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007599"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.565"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'6.404"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.548"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'169.42"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.81"
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.093"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.390"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1008"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.626"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.534"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.338"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04844"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7378"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.132"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.725"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01963"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.879"
$ws.Range("E39").Value = "  +8.23%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.634"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'76.32"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.038"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8732"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4363"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'105.91"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.038.50"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.563"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.249"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1209"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.26"
$ws.Range("E51").Value = "  +0.05%  "
